$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be stored as text so values like "1.026" or
# "27.579.70" are not auto-converted to numbers by Excel's smart typing.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.579.70"
$ws.Range("E2").Value = "  +4.58%  "
$ws.Range("D3").Value = "1.844.52"
$ws.Range("E3").Value = "  +3.99%  "
$ws.Range("D4").Value = "1.026"
$ws.Range("E4").Value = "  +2.31%  "
$ws.Range("D5").Value = "318.75"
$ws.Range("E5").Value = "  +4.11%  "
$ws.Range("D6").Value = "1.025"
$ws.Range("E6").Value = "  +2.34%  "
$ws.Range("D7").Value = "0.4382"
$ws.Range("E7").Value = "  +3.57%  "
$ws.Range("D8").Value = "0.3741"
$ws.Range("E8").Value = "  +3.67%  "
$ws.Range("D9").Value = "0.07397"
$ws.Range("E9").Value = "  +3.68%  "
$ws.Range("D10").Value = "0.8784"
$ws.Range("E10").Value = "  +4.89%  "
$ws.Range("D11").Value = "21.56"
$ws.Range("E11").Value = "  +5.72%  "
$ws.Range("D12").Value = "1.867.22"
$ws.Range("E12").Value = "  +5.14%  "
$ws.Range("D13").Value = "5.503"
$ws.Range("E13").Value = "  +4.89%  "
$ws.Range("D14").Value = "6.697"
$ws.Range("E14").Value = "  +3.87%  "
$ws.Range("D15").Value = "0.07139"
$ws.Range("E15").Value = "  +3.53%  "
$ws.Range("D16").Value = "82.77"
$ws.Range("E16").Value = "  +4.88%  "
$ws.Range("D17").Value = "1.027"
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("D18").Value = "0.000009022"
$ws.Range("E18").Value = "  +4.29%  "
$ws.Range("D19").Value = "1.024"
$ws.Range("E19").Value = "  +2.26%  "
$ws.Range("D20").Value = "15.43"
$ws.Range("E20").Value = "  +3.54%  "
$ws.Range("D21").Value = "27.585.56"
$ws.Range("E21").Value = "  +4.52%  "
$ws.Range("D22").Value = "5.241"
$ws.Range("E22").Value = "  +2.86%  "
$ws.Range("D23").Value = "11.22"
$ws.Range("E23").Value = "  +2.86%  "
$ws.Range("D24").Value = "2.070.11"
$ws.Range("E24").Value = "  +3.82%  "
$ws.Range("D25").Value = "156.95"
$ws.Range("E25").Value = "  +3.47%  "
$ws.Range("D26").Value = "1.919"
$ws.Range("E26").Value = "  +6.84%  "
$ws.Range("D27").Value = "18.72"
$ws.Range("E27").Value = "  +4.11%  "
$ws.Range("D28").Value = "5.270"
$ws.Range("E28").Value = "  +4.18%  "
$ws.Range("D29").Value = "1.943"
$ws.Range("E29").Value = "  +5.63%  "
$ws.Range("D30").Value = "116.39"
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("D31").Value = "0.09078"
$ws.Range("E31").Value = "  +2.88%  "
$ws.Range("D32").Value = "1.214"
$ws.Range("E32").Value = "  +8.63%  "
$ws.Range("D33").Value = "0.7686"
$ws.Range("E33").Value = "  +5.83%  "
$ws.Range("D34").Value = "4.503"
$ws.Range("E34").Value = "  +4.48%  "
$ws.Range("E35").Value = "  +5.01%  "
$ws.Range("E36").Value = "  +2.64%  "
$ws.Range("D37").Value = "1.148"
$ws.Range("E37").Value = "  +6.13%  "
$ws.Range("D38").Value = "0.01974"
$ws.Range("E38").Value = "  +4.82%  "
$ws.Range("D39").Value = "0.05270"
$ws.Range("E39").Value = "  +3.24%  "
$ws.Range("D40").Value = "0.5191"
$ws.Range("E40").Value = "  +5.56%  "
$ws.Range("D41").Value = "2.801"
$ws.Range("E41").Value = "  +7.72%  "
$ws.Range("D42").Value = "0.1669"
$ws.Range("E42").Value = "  +3.72%  "
$ws.Range("D43").Value = "6.638"
$ws.Range("E43").Value = "  +4.96%  "
$ws.Range("D44").Value = "8.571"
$ws.Range("E44").Value = "  +6.50%  "
$ws.Range("D45").Value = "109.36"
$ws.Range("E45").Value = "  +4.52%  "
$ws.Range("D46").Value = "10.57"
$ws.Range("E46").Value = "  +3.60%  "
$ws.Range("E47").Value = "  +2.72%  "
$ws.Range("D48").Value = "1.708"
$ws.Range("E48").Value = "  +5.16%  "
$ws.Range("D49").Value = "0.4658"
$ws.Range("E49").Value = "  +4.95%  "
$ws.Range("D50").Value = "1.912"
$ws.Range("E50").Value = "  +12.03%  "
$ws.Range("D51").Value = "0.06335"
$ws.Range("E51").Value = "  +2.66%  "

# Restore the default cell style so formatting matches the original sheet
# (only the text content changed, not the look of the cells).
$ws.Range("D2:D51").Style = "Normal"
